$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to be bumped
# from 45280 to 45281 for every data row (rows 2 through 27).
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45280) {
        $cell.Value2 = 45281
    }
}
